$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item(3)

# First, extend the header-row style (s=1) and data-row style (s=2) to the
# new trailing columns (H:N) by copying formats from existing styled cells.
$ws3.Range("G1").Copy()
$ws3.Range("H1:N1").PasteSpecial(-4122)
$ws3.Range("G2").Copy()
$ws3.Range("H2:N2").PasteSpecial(-4122)

# Row 1 becomes the real header row (name, capacity, owner, register_date,
# register_reason, acquire_value, property_category, category, date,
# legislator_name, legislator_id, source_file, index).
$ws3.Range("B1").Value = "name"
$ws3.Range("C1").Value = "capacity"
$ws3.Range("D1").Value = "owner"
$ws3.Range("E1").Value = "register_date"
$ws3.Range("F1").Value = "register_reason"
$ws3.Range("G1").Value = "acquire_value"
$ws3.Range("H1").Value = "property_category"
$ws3.Range("I1").Value = "category"
$ws3.Range("J1").Value = "date"
$ws3.Range("K1").Value = "legislator_name"
$ws3.Range("L1").Value = "legislator_id"
$ws3.Range("M1").Value = "source_file"
$ws3.Range("N1").Value = "index"

# Row 2 keeps the original car record and gains the common trailer columns.
$ws3.Range("A2").Value = 38
$ws3.Range("B2").Value = "納智捷G91SPCA"
$ws3.Range("C2").Value = 2198
$ws3.Range("D2").Value = "陳端梅"
$ws3.Range("E2").Value = "100年04月28H"
$ws3.Range("F2").Value = "買賣"
$ws3.Range("G2").Value = 1012000
$ws3.Range("H2").Value = "land"
$ws3.Range("I2").Value = "normal"
$ws3.Range("J2").Value = "2012-03-03"
$ws3.Range("K2").Value = "孫大千"
$ws3.Range("L2").Value = 919
$ws3.Range("M2").Value = "tmpc261"
$ws3.Range("N2").Value = 38

Write-Host "done"
